# "updated figures 5 and 6"
#
# - The bucket label "(2007-2009]" is replaced everywhere by a new,
#   narrower bucket label "[2008-2009]" (rows 3, 8 and 13 on the "2010"
#   sheet). Row 13's label additionally loses its bold/emphasis style.
# - The summary row 27 label is reworded and wrapped onto multiple lines.
# - Figure 5 (row 26) and figure 6 (row 27) formulas are updated to use
#   different weighting factors.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 1, col C was an empty, style-only cell - drop it entirely.
$ws.Range("C1").Clear()

# Rename the "(2007-2009]" bucket to "[2008-2009]" everywhere it appears.
$ws.Range("A3").Value = "[2008-2009]"
$ws.Range("A8").Value = "[2008-2009]"
$ws.Range("A13").Value = "[2008-2009]"

# Row 13's label no longer carries the bold/alternate-font style.
$ws.Range("A13").Style = "Normal"

# Reword the row 27 label and let it wrap across multiple lines.
$ws.Range("A27").Value = "Main AC Units purchased after the start of 2006 but before the end of 2010"
$ws.Range("A27").WrapText = $true
$ws.Rows.Item(27).RowHeight = 51

# Figure 5: weight the "[2006, 2010]" share by 1/3 instead of 1/2.
$ws.Range("B26").Formula = "=B15+B14*(1/3)"

# Figure 6: weight the "[2006, 2010]" share by 2/3, and divide the
# window/wall unit share by 5 instead of 4.
$ws.Range("B27").Formula = "=B14*2/3+B13+B24/5"

# Leave the selection on the newly updated cell.
$ws.Range("B27").Select() | Out-Null
